$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-27 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-28 Sunday", 2) | Out-Null
$d.Content.Find.Execute("271÷9=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "983÷5=196, 3", 2) | Out-Null
$d.Content.Find.Execute("198÷2=99, 0", $true, $false, $false, $false, $false, $true, 1, $false, "399÷2=199, 1", 2) | Out-Null
$d.Content.Find.Execute("164÷3=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "560÷9=62, 2", 2) | Out-Null
$d.Content.Find.Execute("243÷2=121, 1", $true, $false, $false, $false, $false, $true, 1, $false, "917÷8=114, 5", 2) | Out-Null
$d.Content.Find.Execute("731÷3=243, 2", $true, $false, $false, $false, $false, $true, 1, $false, "644÷8=80, 4", 2) | Out-Null
$d.Content.Find.Execute("635÷8=79, 3", $true, $false, $false, $false, $false, $true, 1, $false, "985÷5=197, 0", 2) | Out-Null
$d.Content.Find.Execute("355÷9=39, 4", $true, $false, $false, $false, $false, $true, 1, $false, "500÷4=125, 0", 2) | Out-Null
$d.Content.Find.Execute("565÷7=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "874÷2=437, 0", 2) | Out-Null
$d.Content.Find.Execute("325÷9=36, 1", $true, $false, $false, $false, $false, $true, 1, $false, "934÷2=467, 0", 2) | Out-Null
$d.Content.Find.Execute("148÷3=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "322÷5=64, 2", 2) | Out-Null
$d.Content.Find.Execute("593÷6=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "975÷3=325, 0", 2) | Out-Null
$d.Content.Find.Execute("935÷5=187, 0", $true, $false, $false, $false, $false, $true, 1, $false, "504÷8=63, 0", 2) | Out-Null
$d.Content.Find.Execute("579÷8=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "319÷5=63, 4", 2) | Out-Null
$d.Content.Find.Execute("490÷9=54, 4", $true, $false, $false, $false, $false, $true, 1, $false, "776÷7=110, 6", 2) | Out-Null
$d.Content.Find.Execute("172÷4=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "627÷2=313, 1", 2) | Out-Null
$d.Content.Find.Execute("366÷6=61, 0", $true, $false, $false, $false, $false, $true, 1, $false, "179÷6=29, 5", 2) | Out-Null
$d.Content.Find.Execute("815÷2=407, 1", $true, $false, $false, $false, $false, $true, 1, $false, "757÷9=84, 1", 2) | Out-Null
$d.Content.Find.Execute("188÷7=26, 6", $true, $false, $false, $false, $false, $true, 1, $false, "876÷5=175, 1", 2) | Out-Null
$d.Content.Find.Execute("263÷5=52, 3", $true, $false, $false, $false, $false, $true, 1, $false, "930÷3=310, 0", 2) | Out-Null
$d.Content.Find.Execute("148÷6=24, 4", $true, $false, $false, $false, $false, $true, 1, $false, "417÷2=208, 1", 2) | Out-Null
$d.Content.Find.Execute("217÷3=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "798÷5=159, 3", 2) | Out-Null
$d.Content.Find.Execute("352÷8=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "860÷3=286, 2", 2) | Out-Null
$d.Content.Find.Execute("108÷2=54, 0", $true, $false, $false, $false, $false, $true, 1, $false, "733÷8=91, 5", 2) | Out-Null
$d.Content.Find.Execute("695÷8=86, 7", $true, $false, $false, $false, $false, $true, 1, $false, "280÷3=93, 1", 2) | Out-Null
$d.Content.Find.Execute("297÷8=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "965÷9=107, 2", 2) | Out-Null
